# Update average_county_temperature (K) and recalculated worst_ashp_cop (R) / best_ashp_cop (S)
# values using newer NOAA temperature data, for the affected NAICS 311513 facility rows.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("K5").Value = 12.51681286549706
$ws.Range("R5").Value = 1.782371783972741
$ws.Range("S5").Value = 1.939565227172176
$ws.Range("K6").Value = 12.51681286549706
$ws.Range("R6").Value = 1.782371783972741
$ws.Range("S6").Value = 1.939565227172176
$ws.Range("K7").Value = 15.74228395061728
$ws.Range("R7").Value = 1.837513876759573
$ws.Range("S7").Value = 2.005936573945218
$ws.Range("K8").Value = 15.74228395061728
$ws.Range("R8").Value = 1.837513876759573
$ws.Range("S8").Value = 2.005936573945218
$ws.Range("K9").Value = 1.925925925925943
$ws.Range("R9").Value = 1.62249843161857
$ws.Range("S9").Value = 1.749494516792324
$ws.Range("K10").Value = 1.925925925925943
$ws.Range("R10").Value = 1.62249843161857
$ws.Range("S10").Value = 1.749494516792324
$ws.Range("K11").Value = 1.925925925925943
$ws.Range("K12").Value = 12.66820987654322
$ws.Range("K13").Value = 12.66820987654322
$ws.Range("R13").Value = 1.784885911058073
$ws.Range("S13").Value = 1.942582169301264
$ws.Range("K14").Value = 12.66820987654322
$ws.Range("R14").Value = 1.784885911058073
$ws.Range("S14").Value = 1.942582169301264
$ws.Range("K17").Value = -3.222222222222223
$ws.Range("R17").Value = 1.554711451758341
$ws.Range("S17").Value = 1.669946025515211
$ws.Range("K18").Value = -3.222222222222223
$ws.Range("R18").Value = 1.554711451758341
$ws.Range("S18").Value = 1.669946025515211
$ws.Range("K25").Value = 1.925925925925943
$ws.Range("R25").Value = 1.62249843161857
$ws.Range("S25").Value = 1.749494516792324
$ws.Range("K26").Value = 1.925925925925943
$ws.Range("R26").Value = 1.62249843161857
$ws.Range("S26").Value = 1.749494516792324
$ws.Range("K27").Value = 20.68981481481483
$ws.Range("R27").Value = 1.929056920423291
$ws.Range("S27").Value = 2.117059768804106
$ws.Range("K28").Value = 20.68981481481483
$ws.Range("K29").Value = 20.68981481481483
$ws.Range("R29").Value = 1.929056920423291
$ws.Range("S29").Value = 2.117059768804106
$ws.Range("K30").Value = 14.96875
$ws.Range("R30").Value = 1.8239809580482
$ws.Range("S30").Value = 1.989608681354817
$ws.Range("K31").Value = 14.96875
$ws.Range("R31").Value = 1.8239809580482
$ws.Range("S31").Value = 1.989608681354817
$ws.Range("K32").Value = 1.925925925925943
$ws.Range("R32").Value = 1.62249843161857
$ws.Range("S32").Value = 1.749494516792324
$ws.Range("K33").Value = 1.925925925925943
$ws.Range("R33").Value = 1.62249843161857
$ws.Range("S33").Value = 1.749494516792324
$ws.Range("K34").Value = -3.222222222222223
$ws.Range("R34").Value = 1.554711451758341
$ws.Range("S34").Value = 1.669946025515211
$ws.Range("K35").Value = -3.222222222222223
$ws.Range("K36").Value = -3.222222222222223
$ws.Range("R36").Value = 1.554711451758341
$ws.Range("S36").Value = 1.669946025515211
